$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2
$ws.Range("D2").Value = "41.602.37"
$ws.Range("E2").Value = "  +0.14%  "

# Row 3: update D3, E3
$ws.Range("D3").Value = "2.459.25"
$ws.Range("E3").Value = "  -1.34%  "

# Row 4: update D4, E4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.47%  "

# Row 5: update D5, E5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.64%  "

# Row 6: update D6, E6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.19%  "

# Row 7: update D7, E7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.546"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.15%  "

# Row 8: update E8
$ws.Range("E8").Value = "  +0.31%  "

# Row 9: update E9
$ws.Range("E9").Value = "  +2.56%  "

# Row 10: update D10, E10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.38"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.26%  "

# Row 11: update D11, E11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0811"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.56%  "

# Row 12: update E12
$ws.Range("E12").Value = "  +0.63%  "

# Row 13: update D13, E13
$ws.Range("D13").Value = "2.838.98"
$ws.Range("E13").Value = "  -1.31%  "

# Row 14: update E14
$ws.Range("E14").Value = "  -0.37%  "

# Row 15: update D15, E15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.01%  "

# Row 16: update D16, E16
$ws.Range("D16").Value = "2.473.91"
$ws.Range("E16").Value = "  +0.20%  "

# Row 17: update D17, E17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.776"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.08%  "

# Row 18: update D18
$ws.Range("D18").Value = "41.624.11"

# Row 19: update D19, E19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.33%  "

# Row 20: update E20
$ws.Range("E20").Value = "  +2.01%  "

# Row 21: update D21, E21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.32%  "

# Row 22: update D22, E22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.94%  "

# Row 23: update D23, E23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.94%  "

# Row 24: update E24
$ws.Range("E24").Value = "  +0.24%  "

# Row 25: update B25, C25, D25, E25
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.03%  "

# Row 26: update B26, C26, D26, E26
$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.16%  "

# Row 27: update D27, E27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.30%  "

# Row 28: update E28
$ws.Range("E28").Value = "  +0.79%  "

# Row 29: update D29, E29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.45%  "

# Row 30: update D30, E30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.52%  "

# Row 31: update D31, E31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.82%  "

# Row 32: update E32
$ws.Range("E32").Value = "  +0.56%  "

# Row 33: update E33
$ws.Range("E33").Value = "  +0.05%  "

# Row 34: update D34, E34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0757"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.06%  "

# Row 35: update D35, E35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.66%  "

# Row 36: update D36, E36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.13%  "

# Row 37: update D37, E37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.54%  "

# Row 38: update E38
$ws.Range("E38").Value = "  +1.28%  "

# Row 39: update E39
$ws.Range("E39").Value = "  +1.24%  "

# Row 40: update E40
$ws.Range("E40").Value = "  -2.62%  "

# Row 41: update E41
$ws.Range("E41").Value = "  -4.05%  "

# Row 42: update E42
$ws.Range("E42").Value = "  +0.31%  "

# Row 43: update D43, E43
$ws.Range("D43").Value = "1.975.39"
$ws.Range("E43").Value = "  +1.26%  "

# Row 44: update B44, C44, D44, E44
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0282"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.78%  "

# Row 45: update B45, C45, D45, E45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.78%  "

# Row 46: update D46, E46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.12%  "

# Row 47: update E47
$ws.Range("E47").Value = "  +2.07%  "

# Row 48: update D48, E48
$ws.Range("D48").Value = "2.696.83"
$ws.Range("E48").Value = "  -1.21%  "

# Row 49: update D49, E49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "96.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.49%  "

# Row 50: update D50, E50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.66%  "

# Row 51: update D51, E51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.05%  "
